$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

# New game results for weeks 4-5 of the 2023 season, going into rows 321-352.
# Rows 321-336 already had season/week (columns A/B) filled in, along with an
# (empty but date-formatted) C321; rows 337+ are brand new.
# Each array element is: season, week, date(serial), team1, team2, score1, score2, home_team
$data = @(
    @(2023, 4, 45198, "MIA", "CIN", 15, 27, "CIN"),
    @(2023, 4, 45201, "MIN", "NO", 28, 25, "NO"),
    @(2023, 4, 45201, "WAS", "DAL", 10, 25, "DAL"),
    @(2023, 4, 45201, "LAC", "HOU", 34, 24, "HOU"),
    @(2023, 4, 45201, "TEN", "IND", 24, 17, "IND"),
    @(2023, 4, 45201, "CLE", "ATL", 20, 23, "ATL"),
    @(2023, 4, 45201, "SEA", "DET", 48, 45, "DET"),
    @(2023, 4, 45201, "BUF", "BAL", 23, 20, "BAL"),
    @(2023, 4, 45201, "CHI", "NYG", 12, 20, "NYG"),
    @(2023, 4, 45201, "JAX", "PHI", 21, 29, "PHI"),
    @(2023, 4, 45201, "NYJ", "PIT", 24, 20, "PIT"),
    @(2023, 4, 45201, "ARI", "CAR", 26, 16, "CAR"),
    @(2023, 4, 45201, "DEN", "LV", 23, 32, "LV"),
    @(2023, 4, 45201, "NE", "GB", 24, 27, "GB"),
    @(2023, 4, 45201, "KC", "TB", 41, 31, "TB"),
    @(2023, 4, 45202, "LA", "SF", 9, 24, "SF"),
    @(2023, 5, 45205, "IND", " DEN", 12, 9, "DEN"),
    @(2023, 5, 45208, "NYG", "GB", 27, 22, "GB"),
    @(2023, 5, 45208, "CHI", "MIN", 22, 29, "MIN"),
    @(2023, 5, 45208, "SEA", "NO", 32, 39, "NO"),
    @(2023, 5, 45208, "MIA", "NYJ", 17, 40, "NYJ"),
    @(2023, 5, 45208, "DET", "NE", 0, 29, "NE"),
    @(2023, 5, 45208, "PIT", "BUF", 3, 38, "BUF"),
    @(2023, 5, 45208, "ATL", "TB", 15, 21, "TB"),
    @(2023, 5, 45208, "TEN", "WAS", 21, 17, "WAS"),
    @(2023, 5, 45208, "HOU", "JAX", 13, 6, "JAX"),
    @(2023, 5, 45208, "LAC", "CLE", 30, 28, "CLE"),
    @(2023, 5, 45208, "SF", "CAR", 37, 15, "CAR"),
    @(2023, 5, 45208, "DAL", "LA", 22, 10, "LA"),
    @(2023, 5, 45208, "PHI", "ARI", 20, 17, "ARI"),
    @(2023, 5, 45208, "CIN", "BAL", 17, 19, "BAL"),
    @(2023, 5, 45209, "LV", "KC", 29, 30, "KC")
)

# C320 already carries the date number format used throughout column C;
# reuse it (via Copy) so every new date cell shares that same style instead
# of Excel minting a new one from a NumberFormat string.
$dateFormatSource = $ws.Cells.Item(320, 3)

$startRow = 321
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $dateCell = $ws.Cells.Item($row, 3)
    $dateFormatSource.Copy($dateCell)

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $dateCell.Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}

# Rows 353-373 continue the "season" fill-down pattern that existed before
# (placeholder rows with only the season value filled in).
for ($row = 353; $row -le 373; $row++) {
    $ws.Cells.Item($row, 1).Value = 2023
}

# Leave the selection on the last edited cell, matching the saved workbook state.
$ws.Cells.Item(352, 8).Select()
